$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# ------------------------------------------------------------------
# 1) Insert two new rows at row 44 (pushes nothing since they are at
#    the end of the data, but this registers a proper structural
#    "insert row" edit and inherits row 43's formatting as a start).
# ------------------------------------------------------------------
$ws.Rows.Item(44).Insert()
$ws.Rows.Item(44).Insert()

# ------------------------------------------------------------------
# 2) Re-apply correct cell formatting (borders/fill/number format) by
#    copying formats only from the matching template rows:
#      - row 44 ("DESCE" leg) should look like the plain/unfilled
#        rows (e.g. row 40, the GUARABIRA/SAO PAULO "DESCE" row)
#      - row 45 ("SOBE" leg) should look like the filled/colored
#        rows (e.g. row 41, the GUARABIRA/SAO PAULO "SOBE" row)
# ------------------------------------------------------------------
$ws.Range("A40:M40").Copy() | Out-Null
$ws.Range("A44:M44").PasteSpecial(-4122) | Out-Null

$ws.Range("A41:M41").Copy() | Out-Null
$ws.Range("A45:M45").PasteSpecial(-4122) | Out-Null

# The "SOBE" row's leading (A) cell keeps the plain/unfilled look even
# though the rest of the row is shaded, matching the author's original.
$ws.Range("A40").Copy() | Out-Null
$ws.Range("A45").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3) New shared string used by the two new rows.
# ------------------------------------------------------------------
$linha = "GUARABIRA (PB) - RIO DE JANEIRO (RJ)"

# ------------------------------------------------------------------
# 4) Row 44: GUARABIRA (PB) - RIO DE JANEIRO (RJ) - "DESCE" leg
# ------------------------------------------------------------------
$ws.Range("A44").Formula = "=D44&`" - `"&G44&`" - `"&B44&`" - `"&TEXT(C44,`"HH:MM`")"
$ws.Range("B44").Value = "TER"
$ws.Range("C44").Value = 0.375
$ws.Range("D44").Value = $linha
$ws.Range("E44").Value = "ITAPEMIRIM"
$ws.Range("F44").Value = "IDA"
$ws.Range("G44").Value = "`"DESCE`""
$ws.Range("H44").Value = "GRB"
$ws.Range("I44").Value = "RJO"
$ws.Range("J44").Value = 45888.375
$ws.Range("K44").Formula = "=VLOOKUP(WEEKDAY(J44,1),Planilha2!`$A:`$B,2,0)"
$ws.Range("L44").Value = 45890.458333333336
$ws.Range("M44").Formula = "=VLOOKUP(WEEKDAY(L44,1),Planilha2!`$A:`$B,2,0)"

# ------------------------------------------------------------------
# 5) Row 45: GUARABIRA (PB) - RIO DE JANEIRO (RJ) - "SOBE" leg
# ------------------------------------------------------------------
$ws.Range("A45").Formula = "=D45&`" - `"&G45&`" - `"&B45&`" - `"&TEXT(C45,`"HH:MM`")"
$ws.Range("B45").Value = "SÁB"
$ws.Range("C45").Value = 0.60416666666666663
$ws.Range("D45").Value = $linha
$ws.Range("E45").Value = "ITAPEMIRIM"
$ws.Range("F45").Value = "VOLTA"
$ws.Range("G45").Value = "`"SOBE`""
$ws.Range("H45").Value = "RJO"
$ws.Range("I45").Value = "GRB"
$ws.Range("J45").Value = 45885.604166666664
$ws.Range("K45").Formula = "=VLOOKUP(WEEKDAY(J45,1),Planilha2!`$A:`$B,2,0)"
$ws.Range("L45").Value = 45887.708333333336
$ws.Range("M45").Formula = "=VLOOKUP(WEEKDAY(L45,1),Planilha2!`$A:`$B,2,0)"

# ------------------------------------------------------------------
# 6) Row 42 (CAMPINA GRANDE "DESCE" leg) had an accidental highlighted
#    style; restore it to the plain/unfilled look used by every other
#    "DESCE" row (matches row 25's format).
# ------------------------------------------------------------------
$ws.Range("A25:M25").Copy() | Out-Null
$ws.Range("A42:M42").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 7) Keep the frozen-pane / selection pointed at the new bottom of the
#    sheet, like the author left it.
# ------------------------------------------------------------------
$ws.Range("A46").Select()

Write-Host "edit complete"
